$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of test case data (row 4)
$ws.Range("A4").Value = "CS"
$ws.Range("B4").Value = "CSTest"
$ws.Range("C4").Value = "test"
$ws.Range("D4").Value = "CS TEST"

# Update the active selection to A3 (as reflected in the saved file)
$ws.Range("A3").Select()
